$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns B (TB), C (d2S), D (K), E (IP), G (sum)
# F (Win) column is unchanged by this edit.
$data = @{
    "2" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
    "3" = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=4.358119930609447 }
    "4" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=8.974608811992548 }
    "5" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
    "6" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    "7" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=8.974608811992548 }
    "8" = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=4.358119930609447 }
    "9" = @{ "B"=0.1190320826869504; "C"=0.306821227259698; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=4.457851494814137 }
    "10" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    "11" = @{ "B"=0.2917716402565462; "C"=0.04071648406533734; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=0.9761466351224579 }
    "12" = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=7.143138311642302 }
    "13" = @{ "B"=0.2917716402565462; "C"=0.306821227259698; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=1.242251378316819 }
    "14" = @{ "B"=0.1190320826869504; "C"=0.306821227259698; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=4.457851494814137 }
    "15" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    "16" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    "17" = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=10.19245300693656; "G"=13.45301510845117 }
    "18" = @{ "B"=0.2917716402565462; "C"=0.306821227259698; "D"=3.537761648806719; "E"=0.4942365360607697; "G"=4.630591052383734 }
    "19" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
    "20" = @{ "B"=0.6606524410359556; "C"=0.306821227259698; "D"=0.7527432677738641; "E"=10.19245300693656; "G"=11.91266994300607 }
    "21" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
    "22" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    "23" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=5.586269137925634 }
    "24" = @{ "B"=3.286832544864788; "C"=1.655778082260271; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=6.189590430959694 }
    "25" = @{ "B"=1.455362044514542; "C"=1.655778082260271; "D"=0.1494219747398047; "E"=0.4942365360607697; "G"=3.754798637575387 }
    "26" = @{ "B"=0.1190320826869504; "C"=0.306821227259698; "D"=0.7527432677738641; "E"=0.4942365360607697; "G"=1.672833113781282 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}

Write-Host "done"